$wb = $excel.ActiveWorkbook

# --- Sheet "TC22_Verify_PDP_Page": remove the teststep row that verifies the
#     "Price" text (row 21: VERIFY_TEXT_PRESENT | Price | CSS | Price). All
#     rows below it shift up by one.
$wsTest = $wb.Worksheets.Item("TC22_Verify_PDP_Page")
$wsTest.Rows.Item(21).Delete()

# --- Sheet "Testdata": Baseurl value updated to point at the KIT site.
$wsData = $wb.Worksheets.Item("Testdata")
$wsData.Range("B2").Value = '$BaseURLKIT'
